$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Input@Sink1" row right after each existing data row (rows 2..19),
# working from the bottom up so row indices of rows not yet processed stay stable.
for ($r = 19; $r -ge 2; $r--) {
    $orderId = $ws.Cells.Item($r, 2).Value2

    # Push rows down by inserting a new blank row right after row $r
    $ws.Rows.Item($r + 1).Insert()

    $ws.Cells.Item($r + 1, 1).Value = "Input@Sink1"
    $ws.Cells.Item($r + 1, 2).Value = $orderId
    $ws.Cells.Item($r + 1, 3).Value = 0
    $ws.Cells.Item($r + 1, 4).Value = 0
}
